$p = $ppt.ActivePresentation

$oldDate = "2026. 2. 24."
$newDate = "2026. 2. 27."

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master date placeholder
Update-DateShapes $p.SlideMaster.Shapes

# Every slide layout's date placeholder
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateShapes $layout.Shapes
}

# Notes master date placeholder
Update-DateShapes $p.NotesMaster.Shapes
